# Daily attendance processing - 2026-01-02 19:52:51
# Reorders the comma-separated "Recorded By" values in column G so that the
# previously-trailing entries move to the front (the specific, exact mapping
# mirrors what was observed for each distinct string value in the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact string -> string remap for the "Recorded By" column (column G).
# Only these specific values are touched; any other value (including the
# single-author cells and the untouched 'admin@admin.com, System' combo)
# is left exactly as-is.
$map = @{
    "backup@backdoor.com, system, System" = "system, System, backup@backdoor.com"
    "System, dnasr281@gmail.com"           = "dnasr281@gmail.com, System"
    "backup@backdoor.com, System"          = "System, backup@backdoor.com"
    "admin@admin.com, dnasr281@gmail.com"  = "dnasr281@gmail.com, admin@admin.com"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value2
    if ($map.ContainsKey($val)) {
        $cell.Value2 = $map[$val]
    }
}
